$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item(1)

# C1 -> border: top+bottom thin only (no left/right)
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.LineStyle = 1
$c1.Borders.Item(7).LineStyle = 0
$c1.Borders.Item(10).LineStyle = 0

# D1 -> border: top+right+bottom thin (no left)
$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.LineStyle = 1
$d1.Borders.Item(7).LineStyle = 0

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item(2)

$c1b = $ws2.Range("C1")
$c1b.Style = "Normal"
$c1b.Borders.LineStyle = 1
$c1b.Borders.Item(7).LineStyle = 0
$c1b.Borders.Item(10).LineStyle = 0

$d1b = $ws2.Range("D1")
$d1b.Style = "Normal"
$d1b.Borders.LineStyle = 1
$d1b.Borders.Item(7).LineStyle = 0

$f1b = $ws2.Range("F1")
$f1b.Style = "Normal"
$f1b.Borders.LineStyle = 1
$f1b.Borders.Item(7).LineStyle = 0
$f1b.Borders.Item(10).LineStyle = 0

$g1b = $ws2.Range("G1")
$g1b.Style = "Normal"
$g1b.Borders.LineStyle = 1
$g1b.Borders.Item(7).LineStyle = 0

# Anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5
$ws2.Range("G5").ClearContents()
